$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 3836.0981
$ws.Range("I64").Value = 3629.3076
$ws.Range("J64").Value = 4051.16
$ws.Range("K64").Value = 3629.3076
$ws.Range("L64").Value = 4051.16
$ws.Range("M64").Value = -3381.3076
$ws.Range("N64").Value = -4547.16

# Row 67
$ws.Range("H67").Value = 3836.0981
$ws.Range("I67").Value = 3629.3076
$ws.Range("J67").Value = 4051.16
$ws.Range("K67").Value = 3629.3076
$ws.Range("L67").Value = 4051.16
$ws.Range("M67").Value = -2771.3076
$ws.Range("N67").Value = -5767.16

# Row 76
$ws.Range("H76").Value = 3288.7705
$ws.Range("I76").Value = 2934.5293
$ws.Range("J76").Value = 3425.6365
$ws.Range("K76").Value = 2934.5293
$ws.Range("L76").Value = 3425.6365
$ws.Range("M76").Value = -2619.5293
$ws.Range("N76").Value = -4055.6365

# Row 79
$ws.Range("H79").Value = 3288.7705
$ws.Range("I79").Value = 2934.5293
$ws.Range("J79").Value = 3425.6365
$ws.Range("K79").Value = 2934.5293
$ws.Range("L79").Value = 3425.6365
$ws.Range("M79").Value = -1842.5293
$ws.Range("N79").Value = -5609.636500000001

# Row 98
$ws.Range("H98").Value = 911.3929000000001
$ws.Range("I98").Value = 911.3929000000001
$ws.Range("K98").Value = 911.3929000000001
$ws.Range("M98").Value = 586.6070999999999

# Row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

# Row 122
$ws.Range("H122").Value = 911.3929000000001
$ws.Range("I122").Value = 911.3929000000001
$ws.Range("K122").Value = 2734.1787
$ws.Range("M122").Value = -284.1787000000004

# Row 125
$ws.Range("H125").Value = 2755.5557
$ws.Range("J125").Value = 2755.5557
$ws.Range("L125").Value = 24800.0013
$ws.Range("N125").Value = -29720.0013

# Row 132
$ws.Range("H132").Value = 1773.7561
$ws.Range("I132").Value = 1864.6285
$ws.Range("J132").Value = 1243.6666
$ws.Range("K132").Value = 5593.8855
$ws.Range("L132").Value = 3730.9998
$ws.Range("M132").Value = -3063.8855
$ws.Range("N132").Value = -8790.9998

# Row 138
$ws.Range("H138").Value = 3935.5916
$ws.Range("J138").Value = 4785.7964
$ws.Range("L138").Value = 14357.3892
$ws.Range("N138").Value = -24637.3892

$ws = $wb.Worksheets.Item("ARM")
# Row 88
$ws.Range("H88").Value = 2375
$ws.Range("I88").Value = 2375
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 2375
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -1969
$ws.Range("N88").ClearContents()

# Row 91
$ws.Range("H91").Value = 2375
$ws.Range("I91").Value = 2375
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 2375
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -971
$ws.Range("N91").ClearContents()

# Row 97
$ws.Range("H97").Value = 1262.625
$ws.Range("I97").Value = 1022.2
$ws.Range("J97").Value = 1663.3334
$ws.Range("K97").Value = 1022.2
$ws.Range("L97").Value = 1663.3334
$ws.Range("M97").Value = -526.2
$ws.Range("N97").Value = -2655.3334

# Row 102
$ws.Range("H102").Value = 1801.6
$ws.Range("I102").Value = 1877
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 1877
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = -255
$ws.Range("N102").Value = -4744

$ws = $wb.Worksheets.Item("CRP")
# Row 74
$ws.Range("H74").Value = 22300
$ws.Range("I74").Value = 22800
$ws.Range("J74").Value = 21966.666
$ws.Range("K74").Value = 22800
$ws.Range("L74").Value = 21966.666
$ws.Range("M74").Value = -21926
$ws.Range("N74").Value = -23714.666

# Row 77
$ws.Range("H77").Value = 22300
$ws.Range("I77").Value = 22800
$ws.Range("J77").Value = 21966.666
$ws.Range("K77").Value = 68400
$ws.Range("L77").Value = 65899.99800000001
$ws.Range("M77").Value = -64032
$ws.Range("N77").Value = -74635.99800000001

# Row 106
$ws.Range("H106").Value = 16223.667
$ws.Range("J106").Value = 16223.667
$ws.Range("L106").Value = 16223.667
$ws.Range("N106").Value = -18747.667

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 129.92857
$ws.Range("I23").Value = 97.71429000000001
$ws.Range("J23").Value = 162.14285
$ws.Range("K23").Value = 293.14287
$ws.Range("L23").Value = 486.42855
$ws.Range("M23").Value = -58.14287000000002
$ws.Range("N23").Value = -956.4285500000001

# Row 131
$ws.Range("H131").Value = 864.8182
$ws.Range("I131").Value = 574.75
$ws.Range("J131").Value = 904.8276
$ws.Range("K131").Value = 1724.25
$ws.Range("L131").Value = 2714.4828
$ws.Range("M131").Value = 3315.75
$ws.Range("N131").Value = -12794.4828

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2845.9697
$ws.Range("I80").Value = 2754.0386
$ws.Range("K80").Value = 2754.0386
$ws.Range("M80").Value = -1756.0386

# Row 83
$ws.Range("H83").Value = 2845.9697
$ws.Range("I83").Value = 2754.0386
$ws.Range("K83").Value = 13770.193
$ws.Range("M83").Value = -8778.192999999999

# Row 97
$ws.Range("H97").Value = 1393.1111
$ws.Range("I97").Value = 1317.25
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 1317.25
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -821.25
$ws.Range("N97").Value = -2992

# Row 122
$ws.Range("H122").Value = 7143968
$ws.Range("I122").Value = 7693419.5
$ws.Range("J122").Value = 1100
$ws.Range("K122").Value = 23080258.5
$ws.Range("L122").Value = 3300
$ws.Range("M122").Value = -23077808.5
$ws.Range("N122").Value = -8200

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 9846.412
$ws.Range("I93").Value = 11527.714
$ws.Range("J93").Value = 2000.3334
$ws.Range("K93").Value = 11527.714
$ws.Range("L93").Value = 2000.3334
$ws.Range("M93").Value = -10279.714
$ws.Range("N93").Value = -4496.3334

# Row 100
$ws.Range("H100").Value = 36430076
$ws.Range("I100").Value = 834962.4399999999
$ws.Range("J100").Value = 250000750
$ws.Range("K100").Value = 834962.4399999999
$ws.Range("L100").Value = 250000750
$ws.Range("M100").Value = -834421.4399999999
$ws.Range("N100").Value = -250001832

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 2090
$ws.Range("J81").Value = 2090
$ws.Range("L81").Value = 4180
$ws.Range("N81").Value = -6302

# Row 84
$ws.Range("H84").Value = 2090
$ws.Range("J84").Value = 2090
$ws.Range("L84").Value = 20900
$ws.Range("N84").Value = -31508

# Row 103
$ws.Range("H103").Value = 35000
$ws.Range("J103").Value = 35000
$ws.Range("L103").Value = 35000
$ws.Range("N103").Value = -37344

# Row 122
$ws.Range("H122").Value = 2564.3635
$ws.Range("I122").Value = 2106.8572
$ws.Range("J122").Value = 3365
$ws.Range("K122").Value = 6320.571599999999
$ws.Range("L122").Value = 10095
$ws.Range("M122").Value = -3870.571599999999
$ws.Range("N122").Value = -14995

